# Refresh the per-coin Price (D) and Volume(1h) (E) columns with the
# latest scrape values (GitHub Actions cron update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Column D holds price strings that look numeric (e.g. "606.53",
    # "0.0850", "66.094.42"); Excel auto-coerces bare numeric-looking
    # assignments into floats (dropping trailing zeros / renormalizing)
    # unless the cell is explicitly text-formatted first. Force text,
    # write the literal string, then restore the default/general
    # formatting + Normal style so no stray number format sticks to
    # the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.123.15"
$ws.Range("E2").Value = "  -0.03%  "
Set-TextValue $ws.Range("D3") "3.564.32"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue $ws.Range("D5") "606.53"
$ws.Range("E5").Value = "  +0.73%  "
Set-TextValue $ws.Range("D6") "145.02"
$ws.Range("E6").Value = "  +1.15%  "
Set-TextValue $ws.Range("D7") "3.562.88"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.22%  "
Set-TextValue $ws.Range("D9") "0.490"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("E12").Value = "  -0.12%  "
Set-TextValue $ws.Range("D13") "4.167.81"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("E15").Value = "  -0.78%  "
Set-TextValue $ws.Range("D16") "3.562.22"
$ws.Range("E16").Value = "  +1.61%  "
Set-TextValue $ws.Range("D17") "66.208.52"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  -0.80%  "
Set-TextValue $ws.Range("D19") "11.41"
$ws.Range("E19").Value = "  +8.70%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  +0.36%  "
Set-TextValue $ws.Range("D22") "429.50"
$ws.Range("E22").Value = "  +2.19%  "
Set-TextValue $ws.Range("D23") "0.612"
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("E24").Value = "  +1.64%  "
Set-TextValue $ws.Range("D25") "3.705.20"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  +1.83%  "
Set-TextValue $ws.Range("D29") "7.96"
$ws.Range("E29").Value = "  -0.36%  "
Set-TextValue $ws.Range("D30") "9.11"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("E31").Value = "  -0.19%  "
Set-TextValue $ws.Range("D32") "25.62"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  -1.46%  "
Set-TextValue $ws.Range("D34") "3.559.96"
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +1.70%  "
Set-TextValue $ws.Range("D38") "7.86"
$ws.Range("E38").Value = "  +2.75%  "
Set-TextValue $ws.Range("D39") "5.61"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  -0.36%  "
Set-TextValue $ws.Range("D41") "175.90"
Set-TextValue $ws.Range("D42") "0.0850"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("E44").Value = "  -0.23%  "
Set-TextValue $ws.Range("D45") "1.94"
$ws.Range("E45").Value = "  +1.44%  "
Set-TextValue $ws.Range("D46") "46.05"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("E48").Value = "  -1.42%  "
Set-TextValue $ws.Range("D49") "23.54"
$ws.Range("E49").Value = "  +9.06%  "
$ws.Range("E50").Value = "  +0.18%  "
Set-TextValue $ws.Range("D51") "2.34"
$ws.Range("E51").Value = "  +0.31%  "
